# Updated cryptos list on Thu Nov 21 20:28:57 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell while forcing text storage
# (so numeric-looking strings like '256.34' or '1.00' are not
# auto-converted to numbers by Excel's normal input parsing),
# then clear the temporary text-number-format so no stray style
# is left attached to the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "98.512.95"
Set-TextValue $ws.Range("E2") "  +4.53%  "

Set-TextValue $ws.Range("D3") "3.359.47"
Set-TextValue $ws.Range("E3") "  +9.22%  "

Set-TextValue $ws.Range("E4") "  +0.03%  "

Set-TextValue $ws.Range("D5") "256.34"
Set-TextValue $ws.Range("E5") "  +9.31%  "

Set-TextValue $ws.Range("D6") "623.33"
Set-TextValue $ws.Range("E6") "  +2.59%  "

Set-TextValue $ws.Range("D7") "1.19"
Set-TextValue $ws.Range("E7") "  +8.60%  "

Set-TextValue $ws.Range("D8") "0.389"
Set-TextValue $ws.Range("E8") "  +3.38%  "

Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  +0.02%  "

Set-TextValue $ws.Range("D10") "3.356.18"
Set-TextValue $ws.Range("E10") "  +9.24%  "

Set-TextValue $ws.Range("D11") "0.797"
Set-TextValue $ws.Range("E11") "  -0.82%  "

Set-TextValue $ws.Range("E12") "  +1.97%  "

Set-TextValue $ws.Range("D13") "98.210.39"
Set-TextValue $ws.Range("E13") "  +4.62%  "

Set-TextValue $ws.Range("D14") "36.06"
Set-TextValue $ws.Range("E14") "  +7.40%  "

Set-TextValue $ws.Range("E15") "  +3.42%  "

Set-TextValue $ws.Range("D16") "3.980.78"
Set-TextValue $ws.Range("E16") "  +9.25%  "

Set-TextValue $ws.Range("E17") "  +3.56%  "

Set-TextValue $ws.Range("D18") "3.359.89"
Set-TextValue $ws.Range("E18") "  +10.21%  "

Set-TextValue $ws.Range("E19") "  +2.51%  "

Set-TextValue $ws.Range("E20") "  +4.21%  "

Set-TextValue $ws.Range("D21") "487.72"
Set-TextValue $ws.Range("E21") "  +11.48%  "

Set-TextValue $ws.Range("D22") "5.87"
Set-TextValue $ws.Range("E22") "  +3.62%  "

Set-TextValue $ws.Range("D23") "0.0000208"
Set-TextValue $ws.Range("E23") "  +10.62%  "

Set-TextValue $ws.Range("D24") "9.15"
Set-TextValue $ws.Range("E24") "  +4.36%  "

Set-TextValue $ws.Range("D25") "5.67"
Set-TextValue $ws.Range("E25") "  +3.52%  "

Set-TextValue $ws.Range("D26") "88.06"
Set-TextValue $ws.Range("E26") "  +4.26%  "

Set-TextValue $ws.Range("E27") "  +2.02%  "

Set-TextValue $ws.Range("D28") "3.537.88"
Set-TextValue $ws.Range("E28") "  +9.70%  "

Set-TextValue $ws.Range("E29") "  +0.03%  "

Set-TextValue $ws.Range("D30") "0.186"
Set-TextValue $ws.Range("E30") "  +5.79%  "

Set-TextValue $ws.Range("D31") "0.243"
Set-TextValue $ws.Range("E31") "  -0.68%  "

Set-TextValue $ws.Range("D32") "0.122"
Set-TextValue $ws.Range("E32") "  -0.96%  "

Set-TextValue $ws.Range("E33") "  +11.85%  "

Set-TextValue $ws.Range("D34") "9.26"
Set-TextValue $ws.Range("E34") "  +3.00%  "

Set-TextValue $ws.Range("D35") "27.25"
Set-TextValue $ws.Range("E35") "  +7.77%  "

# Rows 36/37: swap Bittensor <-> Kaspa
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D36") "0.151"
Set-TextValue $ws.Range("E36") "  -1.93%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D37") "515.77"
Set-TextValue $ws.Range("E37") "  +10.89%  "

Set-TextValue $ws.Range("D38") "7.33"
Set-TextValue $ws.Range("E38") "  -3.99%  "

Set-TextValue $ws.Range("D39") "1.94"
Set-TextValue $ws.Range("E39") "  +4.69%  "

Set-TextValue $ws.Range("D40") "24.88"
Set-TextValue $ws.Range("E40") "  +3.70%  "

Set-TextValue $ws.Range("D41") "0.447"
Set-TextValue $ws.Range("E41") "  +2.48%  "

Set-TextValue $ws.Range("D42") "1.27"
Set-TextValue $ws.Range("E42") "  +1.44%  "

Set-TextValue $ws.Range("D43") "3.65"
Set-TextValue $ws.Range("E43") "  -1.23%  "

Set-TextValue $ws.Range("D44") "3.27"
Set-TextValue $ws.Range("E44") "  +5.79%  "

# Rows 45/46: swap USDe <-> ARBITRUM
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D45") "0.783"
Set-TextValue $ws.Range("E45") "  +17.20%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "1.00"
Set-TextValue $ws.Range("E46") "  +0.01%  "

Set-TextValue $ws.Range("D47") "160.60"
Set-TextValue $ws.Range("E47") "  -0.16%  "

Set-TextValue $ws.Range("E48") "  +6.30%  "

Set-TextValue $ws.Range("E49") "  +7.99%  "

Set-TextValue $ws.Range("D50") "45.46"
Set-TextValue $ws.Range("E50") "  +4.32%  "

Set-TextValue $ws.Range("D51") "4.52"
Set-TextValue $ws.Range("E51") "  +7.06%  "
